$d = $word.ActiveDocument

# Correction of the job year:
#   "Denary Computing Ltd in Bangladesh -2019 (www.denarycomputing.com)"
# becomes
#   "Denary Computing Ltd in Bangladesh (2019 - 2020) (www.denarycomputing.com)"
$enDash = [char]0x2013

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$ok = $find.Execute(
    "Bangladesh -2019",                  # FindText
    $true,                                # MatchCase
    $false,                               # MatchWholeWord
    $false,                               # MatchWildcards
    $false,                               # MatchSoundsLike
    $false,                               # MatchAllWordForms
    $true,                                # Forward
    1,                                    # Wrap (wdFindContinue)
    $false,                               # Format
    "Bangladesh (2019 $enDash 2020)",     # ReplaceWith
    2                                     # Replace (wdReplaceAll)
)

if (-not $ok) {
    throw "Could not find the job-year text to correct"
}
